# CORELIMS-98 - production addition of F3
# Updates the NAME column (column C) values for each shelf block,
# bumping the SLC numbers from SLC34-SLC42 to SLC52-SLC60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Range = "C2:C8";   Value = "SLC52" },
    @{ Range = "C9:C15";  Value = "SLC53" },
    @{ Range = "C16:C22"; Value = "SLC54" },
    @{ Range = "C23:C29"; Value = "SLC55" },
    @{ Range = "C30:C36"; Value = "SLC56" },
    @{ Range = "C37:C43"; Value = "SLC57" },
    @{ Range = "C44:C50"; Value = "SLC58" },
    @{ Range = "C51:C57"; Value = "SLC59" },
    @{ Range = "C58:C64"; Value = "SLC60" }
)

foreach ($update in $updates) {
    $ws.Range($update.Range).Value = $update.Value
}
